# Applies the scheduled-runner market-price refresh to the Leve profit
# columns (H:N) across the affected sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3892.439
$ws.Range("I64").Value = 4166.0713
$ws.Range("K64").Value = 4166.0713
$ws.Range("M64").Value = -3918.0713
$ws.Range("H67").Value = 3892.439
$ws.Range("I67").Value = 4166.0713
$ws.Range("K67").Value = 4166.0713
$ws.Range("M67").Value = -3308.0713
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H76").Value = 21298.715
$ws.Range("I76").Value = 32433.646
$ws.Range("J76").Value = 4090.182
$ws.Range("K76").Value = 32433.646
$ws.Range("L76").Value = 4090.182
$ws.Range("M76").Value = -32118.646
$ws.Range("N76").Value = -4720.182
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H79").Value = 21298.715
$ws.Range("I79").Value = 32433.646
$ws.Range("J79").Value = 4090.182
$ws.Range("K79").Value = 32433.646
$ws.Range("L79").Value = 4090.182
$ws.Range("M79").Value = -31341.646
$ws.Range("N79").Value = -6274.182
$ws.Range("H88").Value = 8927.200000000001
$ws.Range("I88").Value = 1900
$ws.Range("J88").Value = 10684
$ws.Range("K88").Value = 1900
$ws.Range("L88").Value = 10684
$ws.Range("M88").Value = -1494
$ws.Range("N88").Value = -11496
$ws.Range("H91").Value = 8927.200000000001
$ws.Range("I91").Value = 1900
$ws.Range("J91").Value = 10684
$ws.Range("K91").Value = 1900
$ws.Range("L91").Value = 10684
$ws.Range("M91").Value = -496
$ws.Range("N91").Value = -13492
$ws.Range("H112").Value = 28572558
$ws.Range("J112").Value = 33614704
$ws.Range("L112").Value = 100844112
$ws.Range("N112").Value = -100846328
$ws.Range("H116").Value = 10318
$ws.Range("I116").Value = 18832.5
$ws.Range("J116").Value = 3019.8572
$ws.Range("K116").Value = 18832.5
$ws.Range("L116").Value = 3019.8572
$ws.Range("M116").Value = -15390.5
$ws.Range("N116").Value = -9903.8572
$ws.Range("H124").Value = 23333.334
$ws.Range("J124").Value = 23333.334
$ws.Range("L124").Value = 23333.334
$ws.Range("N124").Value = -33153.334
$ws.Range("H125").Value = 4495.778
$ws.Range("I125").Value = 4577
$ws.Range("J125").Value = 4333.3335
$ws.Range("K125").Value = 41193
$ws.Range("L125").Value = 39000.0015
$ws.Range("M125").Value = -38733
$ws.Range("N125").Value = -43920.0015

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2432.6667
$ws.Range("I88").Value = 2065.3333
$ws.Range("K88").Value = 2065.3333
$ws.Range("M88").Value = -1659.3333
$ws.Range("H91").Value = 2432.6667
$ws.Range("I91").Value = 2065.3333
$ws.Range("K91").Value = 2065.3333
$ws.Range("M91").Value = -661.3332999999998
$ws.Range("H122").Value = 1833859.6
$ws.Range("I122").Value = 2139036.2
$ws.Range("K122").Value = 6417108.600000001
$ws.Range("M122").Value = -6414658.600000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 11750
$ws.Range("J76").Value = 11750
$ws.Range("L76").Value = 11750
$ws.Range("N76").Value = -12380
$ws.Range("H79").Value = 11750
$ws.Range("J79").Value = 11750
$ws.Range("L79").Value = 11750
$ws.Range("N79").Value = -13934
$ws.Range("H105").Value = 15052
$ws.Range("I105").Value = 22909.9
$ws.Range("K105").Value = 22909.9
$ws.Range("M105").Value = -21162.9

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2685
$ws.Range("J22").Value = 2782.8572
$ws.Range("L22").Value = 8348.571599999999
$ws.Range("N22").Value = -8686.571599999999
$ws.Range("H27").Value = 2685
$ws.Range("J27").Value = 2782.8572
$ws.Range("L27").Value = 8348.571599999999
$ws.Range("N27").Value = -8552.571599999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6567.222
$ws.Range("I70").Value = 6980.048
$ws.Range("J70").Value = 5122.3335
$ws.Range("K70").Value = 6980.048
$ws.Range("L70").Value = 5122.3335
$ws.Range("M70").Value = -6710.048
$ws.Range("N70").Value = -5662.3335
$ws.Range("H73").Value = 6567.222
$ws.Range("I73").Value = 6980.048
$ws.Range("J73").Value = 5122.3335
$ws.Range("K73").Value = 6980.048
$ws.Range("L73").Value = 5122.3335
$ws.Range("M73").Value = -6044.048
$ws.Range("N73").Value = -6994.3335
$ws.Range("H80").Value = 2502.7273
$ws.Range("I80").Value = 2504.6155
$ws.Range("K80").Value = 2504.6155
$ws.Range("M80").Value = -1506.6155
$ws.Range("H83").Value = 2502.7273
$ws.Range("I83").Value = 2504.6155
$ws.Range("K83").Value = 12523.0775
$ws.Range("M83").Value = -7531.077499999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2021.8889
$ws.Range("I7").Value = 2024.625
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 2024.625
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -1912.625
$ws.Range("N7").Value = -2224
$ws.Range("H22").Value = 744.3333
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 774.875
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 774.875
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -1364.875
$ws.Range("H27").Value = 744.3333
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 774.875
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 774.875
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -988.875
$ws.Range("H46").Value = 1208.1666
$ws.Range("I46").Value = 999.5
$ws.Range("J46").Value = 1249.9
$ws.Range("K46").Value = 999.5
$ws.Range("L46").Value = 1249.9
$ws.Range("M46").Value = -811.5
$ws.Range("N46").Value = -1625.9
$ws.Range("H93").Value = 1475
$ws.Range("I93").Value = 1400
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 1400
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = -152
$ws.Range("N93").Value = -3996
$ws.Range("H126").Value = 2021.8889
$ws.Range("I126").Value = 2024.625
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 6073.875
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3603.875
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 20842824
$ws.Range("I132").Value = 41682480
$ws.Range("K132").Value = 125047440
$ws.Range("M132").Value = -125044910

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 43124.5
$ws.Range("J63").Value = 43124.5
$ws.Range("L63").Value = 43124.5
$ws.Range("N63").Value = -44372.5
$ws.Range("H66").Value = 43124.5
$ws.Range("J66").Value = 43124.5
$ws.Range("L66").Value = 129373.5
$ws.Range("N66").Value = -135613.5
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H132").Value = 2215.9678
$ws.Range("I132").Value = 1504.6
$ws.Range("J132").Value = 3509.3635
$ws.Range("K132").Value = 4513.799999999999
$ws.Range("L132").Value = 10528.0905
$ws.Range("M132").Value = -1983.799999999999
$ws.Range("N132").Value = -15588.0905
